$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 91, pushing the existing rows 91..113 down to 93..115
$ws.Range("A91:A92").EntireRow.Insert()

# New row 91
$ws.Range("A91").Value = 9
$ws.Range("B91").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C91").Value = "Metropolitana"
$ws.Range("D91").Value = 44855
$ws.Range("E91").Value = 13
$ws.Range("F91").Value = 100112005
$ws.Range("G91").Value = "Puerro"
$ws.Range("H91").Value = "Sin especificar"
$ws.Range("I91").Value = "Primera"
$ws.Range("J91").Value = 65
$ws.Range("K91").Value = 9000
$ws.Range("L91").Value = 10000
$ws.Range("M91").Value = 9538
$ws.Range("N91").Value = "$/paquete 20 unidades"
$ws.Range("O91").Value = "Provincia de Melipilla"
$ws.Range("P91").Value = 477
$ws.Range("Q91").Value = 20
$ws.Range("R91").Value = "Hortaliza"

# New row 92
$ws.Range("A92").Value = 9
$ws.Range("B92").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C92").Value = "Metropolitana"
$ws.Range("D92").Value = 44855
$ws.Range("E92").Value = 13
$ws.Range("F92").Value = 100112005
$ws.Range("G92").Value = "Puerro"
$ws.Range("H92").Value = "Sin especificar"
$ws.Range("I92").Value = "Primera"
$ws.Range("J92").Value = 80
$ws.Range("K92").Value = 10000
$ws.Range("L92").Value = 10000
$ws.Range("M92").Value = 10000
$ws.Range("N92").Value = "$/paquete 20 unidades"
$ws.Range("O92").Value = "Provincia de Santiago"
$ws.Range("P92").Value = 500
$ws.Range("Q92").Value = 20
$ws.Range("R92").Value = "Hortaliza"
